$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 - name change, Balance updated, Valeur Calculee / Jours de Stock recomputed
$ws.Range("B17").Value = "JULIETTE SOB KEMDJOU"
$ws.Range("F17").Value = 4595
$ws.Range("G17").Value = -2785.523076923077
$ws.Range("H17").Value = 0.6225845989652599

# Row 63
$ws.Range("F63").Value = 425877
$ws.Range("G63").Value = 280772
$ws.Range("H63").Value = 2.934957444609076

# Row 72
$ws.Range("F72").Value = 223436
$ws.Range("G72").Value = 104568
$ws.Range("H72").Value = 1.879698489080324

# Row 73 - name change
$ws.Range("B73").Value = "BLANDINE PEYEMBOUO"
$ws.Range("F73").Value = 50822
$ws.Range("G73").Value = -87022.96296296295
$ws.Range("H73").Value = 0.3686895691187148

# Row 74 - name change
$ws.Range("B74").Value = "ETS CAMPUS III ETS MOBILE FINANCIAL SERVICES MFS"
$ws.Range("F74").Value = 104
$ws.Range("G74").Value = -4896
$ws.Range("H74").Value = 0.0208

# Row 76
$ws.Range("F76").Value = 1098
$ws.Range("G76").Value = -26943.81818181818
$ws.Range("H76").Value = 0.03915580626337289

# Row 78 - name change
$ws.Range("B78").Value = "LAZARRE BIKEK"
$ws.Range("F78").Value = 104496
$ws.Range("G78").Value = 91696
$ws.Range("H78").Value = 8.16375

# Row 85
$ws.Range("F85").Value = 391865
$ws.Range("G85").Value = 357255
$ws.Range("H85").Value = 11.32230569199653

# Row 87
$ws.Range("F87").Value = 64097
$ws.Range("G87").Value = 39067
$ws.Range("H87").Value = 2.560807031562125

# Row 89
$ws.Range("F89").Value = 3842
$ws.Range("G89").Value = -178327.7
$ws.Range("H89").Value = 0.02109022521308428
